$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 46-49 (data for time=44..47) entirely
$ws.Range("A46:C49").EntireRow.Delete()

# Update B and C values for rows 2-45 with new computed energy consumption data
$ws.Range("B2").Value = 1.423522850589841
$ws.Range("C2").Value = 0.8302617648498666
$ws.Range("B3").Value = 5.968883093905558
$ws.Range("C3").Value = 1.470559796438292
$ws.Range("B4").Value = 6.095961842122042
$ws.Range("C4").Value = 2.111374248392587
$ws.Range("B5").Value = 7.344006957089668
$ws.Range("C5").Value = 3.192062064218167
$ws.Range("B6").Value = 22.31850393937586
$ws.Range("C6").Value = 4.291894715285379
$ws.Range("B7").Value = 22.42980282779203
$ws.Range("C7").Value = 4.864982518383937
$ws.Range("B8").Value = 22.52937834687294
$ws.Range("C8").Value = 5.504374693713213
$ws.Range("B9").Value = 25.44727074964677
$ws.Range("C9").Value = 6.491622296892825
$ws.Range("B10").Value = 36.90561329602991
$ws.Range("C10").Value = 7.10657234029987
$ws.Range("B11").Value = 37.0828934789233
$ws.Range("C11").Value = 8.034179602473252
$ws.Range("B12").Value = 37.33094886443832
$ws.Range("C12").Value = 8.757386027820331
$ws.Range("B13").Value = 40.74527899304967
$ws.Range("C13").Value = 9.782113728920043
$ws.Range("B14").Value = 40.84719295123948
$ws.Range("C14").Value = 10.61780568962529
$ws.Range("B15").Value = 41.62054932008182
$ws.Range("C15").Value = 11.22161393859134
$ws.Range("B16").Value = 45.28556612606147
$ws.Range("C16").Value = 11.9813987082552
$ws.Range("B17").Value = 45.40297664034774
$ws.Range("C17").Value = 12.63190730105172
$ws.Range("B18").Value = 48.61710607391728
$ws.Range("C18").Value = 13.32772715558021
$ws.Range("B19").Value = 48.78063400152516
$ws.Range("C19").Value = 13.96802483331814
$ws.Range("B20").Value = 49.1568393219704
$ws.Range("C20").Value = 14.8799630790671
$ws.Range("B21").Value = 50.39191069112326
$ws.Range("C21").Value = 15.59803809403138
$ws.Range("B22").Value = 53.89421236634761
$ws.Range("C22").Value = 16.22809532180113
$ws.Range("B23").Value = 54.65625629224003
$ws.Range("C23").Value = 16.99099167721097
$ws.Range("B24").Value = 58.34298347274057
$ws.Range("C24").Value = 17.66440783092691
$ws.Range("B25").Value = 60.91248896361743
$ws.Range("C25").Value = 18.43661747137801
$ws.Range("B26").Value = 61.05818940988113
$ws.Range("C26").Value = 19.32416613057641
$ws.Range("B27").Value = 69.87793209286971
$ws.Range("C27").Value = 19.95358270658328
$ws.Range("B28").Value = 70.83488873750169
$ws.Range("C28").Value = 21.04409945206961
$ws.Range("B29").Value = 70.94668724739223
$ws.Range("C29").Value = 21.82968679275908
$ws.Range("B30").Value = 73.52741064037168
$ws.Range("C30").Value = 22.86744070181124
$ws.Range("B31").Value = 73.73533097880473
$ws.Range("C31").Value = 23.70042270763484
$ws.Range("B32").Value = 73.83384499078466
$ws.Range("C32").Value = 24.74223611538526
$ws.Range("B33").Value = 77.69661314831956
$ws.Range("C33").Value = 25.6499191420871
$ws.Range("B34").Value = 80.45412447119324
$ws.Range("C34").Value = 26.25818656747703
$ws.Range("B35").Value = 81.76919904311684
$ws.Range("C35").Value = 27.5776653912243
$ws.Range("B36").Value = 81.89077394226412
$ws.Range("C36").Value = 28.23755829528423
$ws.Range("B37").Value = 83.62350807397453
$ws.Range("C37").Value = 29.01647954120543
$ws.Range("B38").Value = 84.68035233131143
$ws.Range("C38").Value = 29.65464233423145
$ws.Range("B39").Value = 85.07040390259884
$ws.Range("C39").Value = 30.32738643002872
$ws.Range("B40").Value = 85.13876611458504
$ws.Range("C40").Value = 31.20025947355375
$ws.Range("B41").Value = 89.89092407201042
$ws.Range("C41").Value = 32.37393042856549
$ws.Range("B42").Value = 90.00250494270522
$ws.Range("C42").Value = 33.02547521722418
$ws.Range("B43").Value = 92.65043132407891
$ws.Range("C43").Value = 34.00416442283062
$ws.Range("B44").Value = 95.2126389533818
$ws.Range("C44").Value = 34.57169762468471
$ws.Range("B45").Value = 95.37192918045898
$ws.Range("C45").Value = 35.53405837206915

Write-Output "done"